$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (AZNM) - newly added values
$ws.Range("B5").Value = 157.32353080410005
$ws.Range("C5").Value = 740.81715305329351
$ws.Range("D5").Value = 29.186054346327136

# Row 7 (CAMX) - updated values
$ws.Range("B7").Value = 165.55318547868359
$ws.Range("C7").Value = 3447.7849909383463
$ws.Range("D7").Value = 1625.3744538197327

# Row 8 (ERCT) - newly added values
$ws.Range("B8").Value = 303.83950137520037
$ws.Range("C8").Value = 2223.7952341798837
$ws.Range("D8").Value = 543.56818099252303

# Row 11 (MROE) - updated values
$ws.Range("B11").Value = 76.378334042006401
$ws.Range("C11").Value = 549.98600263036258
$ws.Range("D11").Value = 263.51049183335653

# Row 12 (MROW) - updated values
$ws.Range("B12").Value = 1654.8434858016274
$ws.Range("C12").Value = 8890.0093098232683
$ws.Range("D12").Value = 1653.3126085875306

# Row 14 (NEWE) - updated values
$ws.Range("B14").Value = 1660.3703094051693
$ws.Range("C14").Value = 4133.3879897129355
$ws.Range("D14").Value = 625.59930181884943

# Row 17 (NYUP) - updated values
$ws.Range("B17").Value = 16.875003144853633
$ws.Range("C17").Value = 490.30649573926831
$ws.Range("D17").Value = 472.60877742829331

# Row 18 (RFCE) - updated values
$ws.Range("B18").Value = 123.86361584785723
$ws.Range("C18").Value = 2393.1237936601133
$ws.Range("D18").Value = 3271.5438335459153

# Row 19 (RFCM) - updated values
$ws.Range("B19").Value = 78.699637831752852
$ws.Range("C19").Value = 2136.2520789749456
$ws.Range("D19").Value = 1749.0868811285004

# Row 20 (RFCW) - updated values
$ws.Range("B20").Value = 8019.0095121843888
$ws.Range("C20").Value = 28860.401618910164
$ws.Range("D20").Value = 4353.9728981157032

# Row 21 (RMPA) - updated values
$ws.Range("B21").Value = 1267.3610666925517
$ws.Range("C21").Value = 5232.9900787739889
$ws.Range("D21").Value = 238.72791888554764

# Row 22 (SPNO) - updated values
$ws.Range("B22").Value = 339.39462024317106
$ws.Range("C22").Value = 1229.4031026149398
$ws.Range("D22").Value = 364.39491045812395

# Row 23 (SPSO) - updated values
$ws.Range("B23").Value = 538.83970166675851
$ws.Range("C23").Value = 6144.4341751539387
$ws.Range("D23").Value = 1500.7126674152721

# Row 24 (SRMV) - updated values
$ws.Range("B24").Value = 61.717717390886861
$ws.Range("C24").Value = 2636.0001245948847
$ws.Range("D24").Value = 59.738967536685664

# Row 25 (SRMW) - updated values
$ws.Range("B25").Value = 907.58548961933172
$ws.Range("C25").Value = 578.02913586846466
$ws.Range("D25").Value = 2430.0538432839876

# Row 26 (SRSO) - updated values
$ws.Range("B26").Value = 1122.2847966709344
$ws.Range("C26").Value = 7319.4625846461086
$ws.Range("D26").Value = 569.45076489300266

# Row 27 (SRTV) - updated values
$ws.Range("B27").Value = 5666.0801951213871
$ws.Range("C27").Value = 5324.3809191971804
$ws.Range("D27").Value = 1152.2123076851203

# Row 28 (SRVC) - updated values
$ws.Range("B28").Value = 1751.7515395261112
$ws.Range("C28").Value = 6735.1202151051266
$ws.Range("D28").Value = 424.42659291092644
